# "Se muestran los reportes sin guardarlos en el disco"
#
# The task "Validacion de cuit para mostrar mensaje correcto" (row 28) is
# removed entirely, shifting every following row up by one. The task that
# used to need the reports saved to disk ("Borrar reportes del disco",
# now row 32 after the shift) is updated to reflect that the reports are
# shown without being saved: Responsable = Agustina, Progreso = 100%.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole row for "Validacion de cuit para mostrar mensaje correcto".
$ws.Rows.Item(28).Delete()

# "Borrar reportes del disco" is now on row 32 after the deletion above.
# Assign it a responsible person and mark it as fully done (100%), matching
# the progress format used by the surrounding rows.
$ws.Range("B32").Value = "Agustina"
$ws.Range("C32").Value = 1
$ws.Range("C32").NumberFormat = "0%"

# Keep the active selection consistent with where the edit was made.
$ws.Range("C33").Select()
